$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1153.1538
$ws.Range("J17").Value = 1153.1538
$ws.Range("L17").Value = 3459.4614
$ws.Range("N17").Value = -3795.4614
$ws.Range("H86").Value = 4161.1
$ws.Range("J86").Value = 3981
$ws.Range("L86").Value = 3981
$ws.Range("N86").Value = -6227
$ws.Range("H89").Value = 4161.1
$ws.Range("J89").Value = 3981
$ws.Range("L89").Value = 19905
$ws.Range("N89").Value = -31137
$ws.Range("H137").Value = 10640794
$ws.Range("I137").Value = 12501461
$ws.Range("J137").Value = 8413
$ws.Range("K137").Value = 37504383
$ws.Range("L137").Value = 25239
$ws.Range("M137").Value = -37501833
$ws.Range("N137").Value = -30339
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 741804.0600000001
$ws.Range("I32").Value = 868129.8
$ws.Range("K32").Value = 868129.8
$ws.Range("M32").Value = -867842.8
$ws.Range("H70").Value = 34000
$ws.Range("J70").Value = 34000
$ws.Range("L70").Value = 34000
$ws.Range("N70").Value = -34540
$ws.Range("H73").Value = 34000
$ws.Range("J73").Value = 34000
$ws.Range("L73").Value = 34000
$ws.Range("N73").Value = -35872
$ws.Range("H97").Value = 651.925
$ws.Range("I97").Value = 585.3143
$ws.Range("K97").Value = 585.3143
$ws.Range("M97").Value = -89.3143
$ws.Range("H110").Value = 1399.8
$ws.Range("I110").Value = 1399.8
$ws.Range("J110").Value = 0
$ws.Range("K110").Value = 1399.8
$ws.Range("L110").Value = 0
$ws.Range("M110").Value = 645.2
$ws.Range("N110").ClearContents()
$ws.Range("H132").Value = 2724907.8
$ws.Range("I132").Value = 2885137.5
$ws.Range("J132").Value = 1000
$ws.Range("K132").Value = 8655412.5
$ws.Range("L132").Value = 3000
$ws.Range("M132").Value = -8652882.5
$ws.Range("N132").Value = -8060
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 27598.44
$ws.Range("I20").Value = 33348.605
$ws.Range("J20").Value = 3879
$ws.Range("K20").Value = 33348.605
$ws.Range("L20").Value = 3879
$ws.Range("M20").Value = -33101.605
$ws.Range("N20").Value = -4373
$ws.Range("H86").Value = 2792.9443
$ws.Range("I86").Value = 2690.6365
$ws.Range("K86").Value = 2690.6365
$ws.Range("M86").Value = -1567.6365
$ws.Range("H89").Value = 2792.9443
$ws.Range("I89").Value = 2690.6365
$ws.Range("K89").Value = 13453.1825
$ws.Range("M89").Value = -7837.182500000001
$ws.Range("H134").Value = 3473898.2
$ws.Range("I134").Value = 1731.3043
$ws.Range("K134").Value = 5193.9129
$ws.Range("M134").Value = -2658.9129
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1545949.5
$ws.Range("I31").Value = 1590005.2
$ws.Range("J31").Value = 4000
$ws.Range("K31").Value = 1590005.2
$ws.Range("L31").Value = 4000
$ws.Range("M31").Value = -1589710.2
$ws.Range("N31").Value = -4590
$ws.Range("H34").Value = 1545949.5
$ws.Range("I34").Value = 1590005.2
$ws.Range("J34").Value = 4000
$ws.Range("K34").Value = 1590005.2
$ws.Range("L34").Value = 4000
$ws.Range("M34").Value = -1589803.2
$ws.Range("N34").Value = -4404
$ws.Range("H86").Value = 13715.091
$ws.Range("I86").Value = 8723
$ws.Range("J86").Value = 16567.715
$ws.Range("K86").Value = 8723
$ws.Range("L86").Value = 16567.715
$ws.Range("M86").Value = -7600
$ws.Range("N86").Value = -18813.715
$ws.Range("H89").Value = 13715.091
$ws.Range("I89").Value = 8723
$ws.Range("J89").Value = 16567.715
$ws.Range("K89").Value = 43615
$ws.Range("L89").Value = 82838.575
$ws.Range("M89").Value = -37999
$ws.Range("N89").Value = -94070.575
$ws.Range("H94").Value = 4297.625
$ws.Range("I94").Value = 2400
$ws.Range("K94").Value = 2400
$ws.Range("M94").Value = -1949
$ws.Range("H107").Value = 1227.5
$ws.Range("I107").Value = 1257.0588
$ws.Range("K107").Value = 1257.0588
$ws.Range("M107").Value = 662.9412
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 4158171
$ws.Range("I4").Value = 5669547.5
$ws.Range("K4").Value = 17008642.5
$ws.Range("M4").Value = -17008530.5
$ws.Range("H5").Value = 1595084.6
$ws.Range("I5").Value = 1880611.4
$ws.Range("J5").Value = 1293695.2
$ws.Range("K5").Value = 5641834.199999999
$ws.Range("L5").Value = 3881085.6
$ws.Range("M5").Value = -5641722.199999999
$ws.Range("N5").Value = -3881309.6
$ws.Range("H23").Value = 383.42856
$ws.Range("I23").Value = 246.25
$ws.Range("K23").Value = 738.75
$ws.Range("M23").Value = -503.75
$ws.Range("H114").Value = 4989
$ws.Range("I114").Value = 0
$ws.Range("J114").Value = 4989
$ws.Range("K114").Value = 0
$ws.Range("L114").Value = 14967
$ws.Range("M114").ClearContents()
$ws.Range("N114").Value = -21475
$ws.Range("H134").Value = 3445.6667
$ws.Range("I134").Value = 3445.6667
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 10337.0001
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -5267.000100000001
$ws.Range("N134").ClearContents()
$ws.Range("H135").Value = 1595084.6
$ws.Range("I135").Value = 1880611.4
$ws.Range("J135").Value = 1293695.2
$ws.Range("K135").Value = 16925502.6
$ws.Range("L135").Value = 11643256.8
$ws.Range("M135").Value = -16922967.6
$ws.Range("N135").Value = -11648326.8
$ws.Range("H138").Value = 43837.668
$ws.Range("I138").Value = 48067.375
$ws.Range("J138").Value = 10000
$ws.Range("K138").Value = 144202.125
$ws.Range("L138").Value = 30000
$ws.Range("M138").Value = -139062.125
$ws.Range("N138").Value = -40280
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 31365.82
$ws.Range("I70").Value = 29402.436
$ws.Range("K70").Value = 29402.436
$ws.Range("M70").Value = -29132.436
$ws.Range("H73").Value = 31365.82
$ws.Range("I73").Value = 29402.436
$ws.Range("K73").Value = 29402.436
$ws.Range("M73").Value = -28466.436
$ws.Range("H109").Value = 58850
$ws.Range("J109").Value = 58850
$ws.Range("L109").Value = 58850
$ws.Range("N109").Value = -60930
$ws.Range("H126").Value = 5972.3887
$ws.Range("I126").Value = 5868.5835
$ws.Range("J126").Value = 6180
$ws.Range("K126").Value = 17605.7505
$ws.Range("L126").Value = 18540
$ws.Range("M126").Value = -15135.7505
$ws.Range("N126").Value = -23480
$ws.Range("H132").Value = 5454.294
$ws.Range("I132").Value = 5981.567
$ws.Range("J132").Value = 1499.75
$ws.Range("K132").Value = 17944.701
$ws.Range("L132").Value = 4499.25
$ws.Range("M132").Value = -15414.701
$ws.Range("N132").Value = -9559.25
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4433.5557
$ws.Range("I7").Value = 3936.75
$ws.Range("K7").Value = 3936.75
$ws.Range("M7").Value = -3824.75
$ws.Range("H22").Value = 3227.8262
$ws.Range("I22").Value = 1860
$ws.Range("K22").Value = 1860
$ws.Range("M22").Value = -1565
$ws.Range("H27").Value = 3227.8262
$ws.Range("I27").Value = 1860
$ws.Range("K27").Value = 1860
$ws.Range("M27").Value = -1753
$ws.Range("H74").Value = 40877.4
$ws.Range("I74").Value = 40877.4
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 40877.4
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -39879.4
$ws.Range("N74").ClearContents()
$ws.Range("H77").Value = 40877.4
$ws.Range("I77").Value = 40877.4
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 122632.2
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -117640.2
$ws.Range("N77").ClearContents()
$ws.Range("H122").Value = 5857
$ws.Range("I122").Value = 4200
$ws.Range("J122").Value = 9999.5
$ws.Range("K122").Value = 12600
$ws.Range("L122").Value = 29998.5
$ws.Range("M122").Value = -10150
$ws.Range("N122").Value = -34898.5
$ws.Range("H126").Value = 4433.5557
$ws.Range("I126").Value = 3936.75
$ws.Range("K126").Value = 11810.25
$ws.Range("M126").Value = -9340.25
$ws.Range("H132").Value = 1362795.4
$ws.Range("I132").Value = 1962993.6
$ws.Range("J132").Value = 2346.2
$ws.Range("K132").Value = 5888980.800000001
$ws.Range("L132").Value = 7038.599999999999
$ws.Range("M132").Value = -5886450.800000001
$ws.Range("N132").Value = -12098.6
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = 0
$ws.Range("K14").Value = 0
$ws.Range("M14").ClearContents()
$ws.Range("H26").Value = 1012
$ws.Range("I26").Value = 1012
$ws.Range("K26").Value = 1012
$ws.Range("M26").Value = -719
$ws.Range("H45").Value = 34217.5
$ws.Range("J45").Value = 36623.332
$ws.Range("L45").Value = 36623.332
$ws.Range("N45").Value = -37605.332
$ws.Range("H122").Value = 33698.46
$ws.Range("I122").Value = 3103.4
$ws.Range("J122").Value = 97438.164
$ws.Range("K122").Value = 9310.200000000001
$ws.Range("L122").Value = 292314.492
$ws.Range("M122").Value = -6860.200000000001
$ws.Range("N122").Value = -297214.492
$ws.Range("H132").Value = 9261497
$ws.Range("I132").Value = 11906125
$ws.Range("K132").Value = 35718375
$ws.Range("M132").Value = -35715845
$ws.Range("H140").Value = 79594.25
$ws.Range("J140").Value = 79594.25
$ws.Range("L140").Value = 79594.25
$ws.Range("N140").Value = -89954.25
$ws.Range("H141").Value = 81905
$ws.Range("J141").Value = 81905
$ws.Range("L141").Value = 81905
$ws.Range("N141").Value = -92265
